$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 323 (date 44209): new positive cases corrected, one more extra-hospital death ---
$ws.Range("C323").Value = 134

# M323's column is formatted as Text ("@"); temporarily switch to a numeric
# format so the value is stored as a real number (matching the rest of the
# column), then restore the original Text format.
$fmt = $ws.Range("M323").NumberFormat
$ws.Range("M323").NumberFormat = "General"
$ws.Range("M323").Value = 2
$ws.Range("M323").NumberFormat = $fmt

# --- Row 324 (date 44210): new positive cases corrected ---
$ws.Range("C324").Value = 84

# --- Row 325 (date 44211): day was previously blank, now has real data ---
$ws.Range("C325").Value = 19
$ws.Range("E325").Value = 11
$ws.Range("F325").Value = 7
$ws.Range("G325").Value = 124

$fmt = $ws.Range("L325").NumberFormat
$ws.Range("L325").NumberFormat = "General"
$ws.Range("L325").Value = 0
$ws.Range("L325").NumberFormat = $fmt

$fmt = $ws.Range("M325").NumberFormat
$ws.Range("M325").NumberFormat = "General"
$ws.Range("M325").Value = 0
$ws.Range("M325").NumberFormat = $fmt

# --- Update the saved selection of the frozen bottom-right pane from Z2 to A2 ---
$ws.Range("A2").Select()
